$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ============ Row 11 (tall, multi-line row) ============
$ws.Range("A4").Copy()
$ws.Range("A11").PasteSpecial(-4122)

$ws.Range("B4").Copy()
$ws.Range("B11").PasteSpecial(-4122)

$ws.Range("C4").Copy()
$ws.Range("C11").PasteSpecial(-4122)

$ws.Range("D4").Copy()
$ws.Range("D11").PasteSpecial(-4122)

$ws.Range("E4").Copy()
$ws.Range("E11").PasteSpecial(-4122)

$ws.Range("F7").Copy()
$ws.Range("F11").PasteSpecial(-4122)

$ws.Range("G5").Copy()
$ws.Range("G11").PasteSpecial(-4122)

$ws.Range("H4").Copy()
$ws.Range("H11").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false
$ws.Rows.Item(11).RowHeight = 72

$ws.Range("A11").Value = "SRS_02"
$ws.Range("B11").Value = "Software Requirements Specification"
$ws.Range("C11").Value = "Nourhan Ali,Mohamed Ibrahim`n,Manar Ali ,Aalaa Adel and Al-Shimaa`nShehata"
$ws.Range("D11").Value = "Mostafa Mohamed"
$ws.Range("E11").Value = 45508
$ws.Range("E11").NumberFormat = "mm-dd-yy"
$ws.Range("F11").Value = "__"
$ws.Range("G11").Value = "__"
$ws.Range("H11").Value = "Solved"

# ============ Row 12 (normal height row) ============
$ws.Range("A10").Copy()
$ws.Range("A12").PasteSpecial(-4122)

$ws.Range("B5").Copy()
$ws.Range("B12").PasteSpecial(-4122)

$ws.Range("C5").Copy()
$ws.Range("C12").PasteSpecial(-4122)

$ws.Range("D5").Copy()
$ws.Range("D12").PasteSpecial(-4122)

$ws.Range("A10").Copy()
$ws.Range("E12").PasteSpecial(-4122)

$ws.Range("G5").Copy()
$ws.Range("G12").PasteSpecial(-4122)

$ws.Range("H10").Copy()
$ws.Range("H12").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false
$ws.Rows.Item(12).RowHeight = 21

$ws.Range("A12").Value = "RTM_01"
$ws.Range("B12").Value = "Requirment Tracability Matrix"
$ws.Range("C12").Value = "Mostafa Mohamed"
$ws.Range("D12").Value = "Mohamed Ibrahim"
$ws.Range("E12").Value = 45508
$ws.Range("E12").NumberFormat = "mm-dd-yy"

$ws.Range("F12").HorizontalAlignment = -4108
$ws.Range("F12").VerticalAlignment = -4160
$ws.Range("F12").WrapText = $true
$ws.Range("F12").Value = "__"

$ws.Range("G12").Value = "__"
$ws.Range("H12").Value = "Solved"

# ============ sheet view changes ============
$ws.Range("A13").Select()
$excel.ActiveWindow.Zoom = 109
